# Edit: renamed the Spellcaster searcher ("Phantom Sorcerer" -> "Assistant Magician")
# with new art source ("Phantom Dewan" -> "Boo Koo") and ATK/DEF filled in, plus
# rebalanced several searchers' Levels and added a Level 2/3/4 breakdown.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Single-Type Searchers")

# --- Level (column F) rebalances across several rows (searcher name/text unchanged) ---
$ws.Range("F2").Value = 2
$ws.Range("F3").Value = 2
$ws.Range("F9").Value = 3
$ws.Range("F12").Value = 2
$ws.Range("F13").Value = 2
$ws.Range("F14").Value = 2
$ws.Range("F19").Value = 3
$ws.Range("F20").Value = 3
$ws.Range("F24").Value = 2
$ws.Range("F25").Value = 2

# --- Row 11 (Machine / "Machine Soldier"): fill in ATK/DEF, then restyle to match
#     the highlighted look used on row 7 / 15 ---
$ws.Range("H11").Value = 1500
$ws.Range("I11").Value = 1000
$ws.Range("E7:L7").Copy()
$ws.Range("E11:L11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- Row 17 (Spellcaster searcher): finished art, renamed monster + art source ---
$ws.Range("E17").Value = "Assistant Magician"
$ws.Range("F17").Value = 2
$ws.Range("H17").Value = 450
$ws.Range("I17").Value = 750
$ws.Range("L17").Value = "Boo Koo"

# --- New Level 2 / 3 / 4 breakdown next to the existing restriction counts ---
$ws.Range("J27").Value = "Level 2"
$ws.Range("K27").Formula = '=COUNTIF(F2:F25, "=2")'
$ws.Range("J28").Value = "Level 3"
$ws.Range("K28").Formula = '=COUNTIF(F2:F25, "=3")'
$ws.Range("J29").Value = "Level 4"
$ws.Range("K29").Formula = '=COUNTIF(F2:F25, "=4")'

# --- Update the saved selection ---
$ws.Range("J29").Select()
